# Import "ho dem" (surname) as a new column between "mssv" and "ten".
#
# Target layout (Sheet1, row 1 header):
#   A:STT  B:mssv  C:hodem(NEW)  D:ten  E:ngay_sinh  F:lop

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new blank column at C; this shifts ten/ngay_sinh/lop
# (previously C/D/E) one column to the right, becoming D/E/F, and extends
# the used range from A1:E1 to A1:F1.
$ws.Columns.Item(3).Insert()

# Give the header of the freshly inserted column its text.
$ws.Range("C1").Value = "hodem"

# The new column should be the same width as column B (mssv); round to the
# nearest width the engine can represent.
$ws.Columns.Item(3).ColumnWidth = 14.25

# Reflect the selection left behind by the editing session.
[void]$ws.Range("P23").Select()
